$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.025.95"
$ws.Range("E2").Value = "  -1.35%  "
$ws.Range("D3").Value = "1.667.12"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.81"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5107"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2680"
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06374"
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.86"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07435"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "1.676.62"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.505"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008473"
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.02"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("D17").Value = "25.834.82"
$ws.Range("E17").Value = "  -2.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.921"
$ws.Range("E18").Value = "  -1.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.79"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("E21").Value = "  +1.18%  "
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.99"
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.599"
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1215"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.64"
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06623"
$ws.Range("E28").Value = "  +14.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.328"
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.311"
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.508"
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6178"
$ws.Range("E35").Value = "  +3.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.368"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.684"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.294"
$ws.Range("E38").Value = "  +6.71%  "
$ws.Range("D39").Value = "1.095.95"
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8660"
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.52"
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("D44").Value = "1.814.81"
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.21"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.008"
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.117"
$ws.Range("E48").Value = "  +1.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05226"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.980"
$ws.Range("E51").Value = "  +2.39%  "
